$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Cream Cakes sheet: price update on row 12 (500 -> 450)
# ---------------------------------------------------------------------------
$wsCream = $wb.Worksheets.Item("Cream Cakes")
$wsCream.Range("D12").Value = 450

# ---------------------------------------------------------------------------
# Sweet Delights sheet: price update on row 3 (45 -> 35) and five new
# "cookie" products appended as rows 6-10.
# ---------------------------------------------------------------------------
$wsSweet = $wb.Worksheets.Item("Sweet Delights")
$wsSweet.Range("D3").Value = 35

# Pre-format the cells that end up carrying the "data row" style (s="2")
# by copying formats from an existing formatted row before filling values in,
# mirroring how the new rows were styled only partially.
$wsSweet.Range("A2").Copy() | Out-Null
$wsSweet.Range("A6").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("B6").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("D6").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("E6").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("G6").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("H6").PasteSpecial(-4122) | Out-Null

$wsSweet.Range("A7").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("B7").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("D7").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("G7").PasteSpecial(-4122) | Out-Null

$wsSweet.Range("A8").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("B8").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("D8").PasteSpecial(-4122) | Out-Null

$wsSweet.Range("A9").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("B9").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("D9").PasteSpecial(-4122) | Out-Null

$wsSweet.Range("A10").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("B10").PasteSpecial(-4122) | Out-Null
$wsSweet.Range("D10").PasteSpecial(-4122) | Out-Null

# Fill in the values in the order that reproduces the shared-string table
# layout: ids, then names, then the remaining images (jeera's image last).
$wsSweet.Range("A6").Value = "sd5"
$wsSweet.Range("B6").Value = "jeera cookies"
$wsSweet.Range("E6").Value = "For 100 gm"

$wsSweet.Range("A7").Value = "sd6"
$wsSweet.Range("A8").Value = "sd7"
$wsSweet.Range("A9").Value = "sd8"

$wsSweet.Range("B7").Value = "atta cookies"
$wsSweet.Range("B8").Value = "chocolate cookies"
$wsSweet.Range("B9").Value = "coconut cookies"

$wsSweet.Range("A10").Value = "sd9"
$wsSweet.Range("B10").Value = "ragi oats cookies"

$wsSweet.Range("C7").Value = "sweet-delights/atta-cookies.jpg"
$wsSweet.Range("C8").Value = "sweet-delights/chocolate-cookies.jpg"
$wsSweet.Range("C9").Value = "sweet-delights/coconut-cookies.jpg"
$wsSweet.Range("C10").Value = "sweet-delights/ragi-oats-cookies.jpg"
$wsSweet.Range("C6").Value = "sweet-delights/jeera-cookies.jpg"

# Prices
$wsSweet.Range("D6").Value = 30
$wsSweet.Range("D7").Value = 30
$wsSweet.Range("D8").Value = 40
$wsSweet.Range("D9").Value = 35
$wsSweet.Range("D10").Value = 40

# priceFor (re-uses the "For 100 gm" shared string created above)
$wsSweet.Range("E7").Value = "For 100 gm"
$wsSweet.Range("E8").Value = "For 100 gm"
$wsSweet.Range("E9").Value = "For 100 gm"
$wsSweet.Range("E10").Value = "For 100 gm"

# inStock / onDiscount (re-use the existing "yes"/"no" shared strings)
$wsSweet.Range("G6").Value = "yes"
$wsSweet.Range("G7").Value = "yes"
$wsSweet.Range("G8").Value = "yes"
$wsSweet.Range("G9").Value = "yes"
$wsSweet.Range("G10").Value = "yes"

$wsSweet.Range("H6").Value = "no"
$wsSweet.Range("H7").Value = "no"
$wsSweet.Range("H8").Value = "no"
$wsSweet.Range("H9").Value = "no"
$wsSweet.Range("H10").Value = "no"

# Row heights to match the rest of the table (15.75pt, custom height)
$wsSweet.Rows.Item(6).RowHeight = 15.75
$wsSweet.Rows.Item(7).RowHeight = 15.75
$wsSweet.Rows.Item(8).RowHeight = 15.75
$wsSweet.Rows.Item(9).RowHeight = 15.75
$wsSweet.Rows.Item(10).RowHeight = 15.75

# ---------------------------------------------------------------------------
# Selections / active sheet: Pastries was the selected tab before, now it's
# Sweet Delights (matches activeTab 2 -> 3 and the per-sheet tabSelected flag
# moving from Pastries to Sweet Delights). Cream Cakes' selection also moved.
# ---------------------------------------------------------------------------
$wsCream.Range("D13").Select()

$wsSweet.Activate()
$wsSweet.Range("D11").Select()

Write-Host "done"
